$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 596, pushing the existing rows 596:626 down to 597:627
$ws.Rows(596).Insert()

# Populate the newly inserted row 596 with the new weekly record
$ws.Cells.Item(596, 1).Value = 3
$ws.Cells.Item(596, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(596, 3).Value = "Coquimbo"
$ws.Cells.Item(596, 4).Value = 44585
$ws.Cells.Item(596, 5).Value = 5
$ws.Cells.Item(596, 6).Value = "Fruta"
$ws.Cells.Item(596, 7).Value = 100102
$ws.Cells.Item(596, 8).Value = "Cítricos"
$ws.Cells.Item(596, 9).Value = 100102005
$ws.Cells.Item(596, 10).Value = "Naranja"
$ws.Cells.Item(596, 11).Value = "Valencia"
$ws.Cells.Item(596, 12).Value = "Primera"
$ws.Cells.Item(596, 13).Value = 170
$ws.Cells.Item(596, 14).Value = 6000
$ws.Cells.Item(596, 15).Value = 6500
$ws.Cells.Item(596, 16).Value = 6235
$ws.Cells.Item(596, 17).Value = "$/malla 13 kilos"
$ws.Cells.Item(596, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(596, 19).Value = 480
$ws.Cells.Item(596, 20).Value = 13
